# Replace "OIE" with "WOAH" throughout the LSDV story-map content,
# matching the author's commit: "OIE replaced with WOAH all Excels".
#
# The affected cells live on two worksheets:
#   "Sheet 1"    -> column E ("Content") for several rows
#   "References" -> column C ("Paper") for several rows
#
# Hyperlink URLs (which still contain lowercase "oie.int") are left
# untouched, matching the diff.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$sheet2 = $wb.Worksheets.Item("References")

function Replace-OIE {
    param(
        [object]$ws,
        [string]$addr
    )
    $cell = $ws.Range($addr)
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current.Replace("OIE", "WOAH")
    }
}

# Sheet 1 ("Content" column E)
$sheet1Addresses = @("E5", "E6", "E7", "E14", "E17", "E21", "E34", "E43", "E55", "E68", "E94")
foreach ($addr in $sheet1Addresses) {
    Replace-OIE -ws $sheet1 -addr $addr
}

# References sheet ("Paper" column C)
$sheet2Addresses = @("C2", "C5", "C8", "C9", "C10")
foreach ($addr in $sheet2Addresses) {
    Replace-OIE -ws $sheet2 -addr $addr
}
